# Add two new columns I (I0) and J (IF) to the worksheet, mirroring the
# existing header style used by columns A-H, and populate the data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers: I1 = "I0", J1 = "IF", using the same style as the other header cells.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-37 for columns I and J.
$values = @{
    2  = @(1,4)
    3  = @(1,6)
    4  = @(1,5)
    5  = @(1,5)
    6  = @(1,6)
    7  = @(1,6)
    8  = @(1,9)
    9  = @(1,5)
    10 = @(1,6)
    11 = @(1,3)
    12 = @(1,6)
    13 = @(1,6)
    14 = @(1,5)
    15 = @(1,4)
    16 = @(1,6)
    17 = @(1,6)
    18 = @(1,4)
    19 = @(1,5)
    20 = @(1,5)
    21 = @(9,9)
    22 = @(3,4)
    23 = @(10,11)
    24 = @(6,6)
    25 = @(4,6)
    26 = @(7,8)
    27 = @(6,8)
    28 = @(5,6)
    29 = @(1,4)
    30 = @(2,4)
    31 = @(7,8)
    32 = @(9,9)
    33 = @(9,9)
    34 = @(1,5)
    35 = @(1,3)
    36 = @(1,3)
    37 = @(1,2)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
